# Auto-generated edit script: update crypto price (D) and volume-change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.844.12"
$c.ClearFormats()
$ws.Range("E2").Value = "  -1.76%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.824.16"
$c.ClearFormats()
$ws.Range("E3").Value = "  -1.71%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.ClearFormats()
$ws.Range("E4").Value = "  +0.57%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "310.26"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.12%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.ClearFormats()
$ws.Range("E6").Value = "  +0.46%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4566"
$c.ClearFormats()
$ws.Range("E7").Value = "  -0.86%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3674"
$c.ClearFormats()
$ws.Range("E8").Value = "  -0.99%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07149"
$c.ClearFormats()
$ws.Range("E9").Value = "  -2.34%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8714"
$c.ClearFormats()
$ws.Range("E10").Value = "  -0.97%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07772"
$c.ClearFormats()
$ws.Range("E11").Value = "  -0.42%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "19.50"
$c.ClearFormats()
$ws.Range("E12").Value = "  -1.87%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.818.30"
$c.ClearFormats()
$ws.Range("E13").Value = "  -1.36%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.308"
$c.ClearFormats()
$ws.Range("E14").Value = "  -1.50%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.370"
$c.ClearFormats()
$ws.Range("E15").Value = "  -2.62%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "86.66"
$c.ClearFormats()
$ws.Range("E16").Value = "  -5.61%  "

$ws.Range("E17").Value = "  +0.54%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008687"
$c.ClearFormats()
$ws.Range("E18").Value = "  -3.89%  "

$ws.Range("E19").Value = "  +0.47%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "26.855.53"
$c.ClearFormats()
$ws.Range("E20").Value = "  -1.79%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.41"
$c.ClearFormats()
$ws.Range("E21").Value = "  -2.50%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.986"
$c.ClearFormats()
$ws.Range("E22").Value = "  -2.79%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "2.047.46"
$c.ClearFormats()
$ws.Range("E23").Value = "  -3.29%  "

$ws.Range("E24").Value = "  -0.91%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.998"
$c.ClearFormats()
$ws.Range("E25").Value = "  +4.13%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "150.95"
$c.ClearFormats()
$ws.Range("E26").Value = "  -0.76%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.16"
$c.ClearFormats()
$ws.Range("E27").Value = "  -1.21%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.945"
$c.ClearFormats()
$ws.Range("E28").Value = "  -6.08%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "113.50"
$c.ClearFormats()
$ws.Range("E29").Value = "  -2.32%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.896"
$c.ClearFormats()
$ws.Range("E30").Value = "  -4.13%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08784"
$c.ClearFormats()
$ws.Range("E31").Value = "  -0.90%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.020"
$c.ClearFormats()
$ws.Range("E32").Value = "  -0.38%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.7454"
$c.ClearFormats()
$ws.Range("E33").Value = "  -3.67%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.467"
$c.ClearFormats()
$ws.Range("E34").Value = "  -0.45%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.127"
$c.ClearFormats()
$ws.Range("E35").Value = "  -4.24%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.535"
$c.ClearFormats()
$ws.Range("E36").Value = "  -3.70%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.084"
$c.ClearFormats()
$ws.Range("E37").Value = "  +0.66%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01932"
$c.ClearFormats()
$ws.Range("E38").Value = "  -1.44%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.923"
$c.ClearFormats()
$ws.Range("E39").Value = "  -0.99%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.05104"
$c.ClearFormats()
$ws.Range("E40").Value = "  -2.41%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.913"
$c.ClearFormats()
$ws.Range("E41").Value = "  -1.68%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.4946"
$c.ClearFormats()
$ws.Range("E42").Value = "  -3.80%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.1591"
$c.ClearFormats()
$ws.Range("E43").Value = "  -2.85%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.260"
$c.ClearFormats()

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.4664"
$c.ClearFormats()
$ws.Range("E45").Value = "  -3.49%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.ClearFormats()
$ws.Range("E46").Value = "  +0.50%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "10.05"
$c.ClearFormats()
$ws.Range("E47").Value = "  -2.78%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "101.07"
$c.ClearFormats()
$ws.Range("E48").Value = "  -2.07%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.605"
$c.ClearFormats()
$ws.Range("E49").Value = "  -2.87%  "

$ws.Range("E50").Value = "  -2.08%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "64.30"
$c.ClearFormats()
$ws.Range("E51").Value = "  -2.38%  "

